$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'40.273.06"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +3.33%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.245.80"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.99%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.05%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'297.17"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -0.16%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'87.33"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +8.20%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +1.94%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.03%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.476"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +3.40%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'31.32"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +11.57%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.0802"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +3.63%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'47.25"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +3.70%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.99%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'6.50"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +6.02%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'2.593.49"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +1.10%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'14.28"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +1.78%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'2.228.38"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -2.21%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'0.736"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +3.05%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'40.193.86"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +3.37%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +4.33%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'5.86"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.91%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'10.89"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +10.07%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'65.75"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.41%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'236.90"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +5.21%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  +0.05%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'2.49"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +3.55%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  +6.79%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'23.05"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +3.54%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +2.53%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'9.30"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +4.39%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'33.69"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +7.42%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +3.68%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.02%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'4.92"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +2.94%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  +4.96%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  +3.19%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'16.74"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +14.95%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  +7.30%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +3.23%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  +1.51%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'1.71"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +7.21%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'3.85"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +6.89%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'2.042.03"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +6.77%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'2.22"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +9.91%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.0273"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +6.80%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'10.06"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +12.67%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'16.35"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.27%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'2.59"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +2.68%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'2.464.12"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +1.61%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'71.68"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +4.29%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'Stacks"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'1.47"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +15.40%  "
$ws.Range('E51').Style = 'Normal'
